# CPA.xlsx update:
#  - add a "Rang" column (F) header
#  - fill in the "Moyenne de l'étudiant" column (E) with each student's average
#  - default font size moves from 12pt to 11pt (workbook "Normal" style)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Default/Normal style font size: 12 -> 11
$wb.Styles.Item("Normal").Font.Size = 11

# New column header
$ws.Range("F1").Value = "Rang"

# Per-student average ("Moyenne de l'étudiant"), rows 3..63 (row 2 stays blank,
# matching the existing gap in the sheet)
$moyennes = @(
    6, 16, 13, 17, 20, 13, 18, 8, 7, 5,
    19, 11, 17, 6, 17, 6, 6, 11, 10, 17,
    11, 10, 9, 5, 9, 10, 7, 14, 7, 19,
    15, 20, 11, 20, 16, 13, 11, 7, 19, 15,
    11, 7, 13, 8, 14, 20, 7, 5, 15, 13,
    9, 11, 6, 13, 9, 12, 13, 13, 10, 15,
    8
)

$firstRow = 3
for ($i = 0; $i -lt $moyennes.Length; $i++) {
    $row = $firstRow + $i
    $ws.Range("E$row").Value = $moyennes[$i]
}
